$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 6.2205449478915495
    "C2" = 10.088932515975207
    "D2" = 7.7910835671158161
    "E2" = 4.2627860862114568

    "B3" = 0.90468654188861353
    "C3" = 1.3002718982054766
    "D3" = 0.91863158209723739
    "E3" = 1.8260152612626817

    "B4" = 4.6341749033787023
    "C4" = 5.0657786490690713
    "D4" = 6.4329677061929003
    "E4" = 6.9795806742906992

    "B5" = 2.5079867320948721
    "C5" = 2.2390278274351099
    "D5" = 2.3199095589483352
    "E5" = 2.0348548335227075

    "E6" = 3.1003992135261251

    "B9" = 25079.867320948721
    "C9" = 22390.278274351098
    "D9" = 23199.09558948335
    "E9" = 20348.548335227075

    "E10" = 31003.992135261251

    "B11" = 4.3342674124273604
    "C11" = 6.1022571540750352
    "D11" = 5.1161769500034797
    "E11" = 3.0215229534704662

    "B12" = 25.850573792598141
    "C12" = 29.97628653383693
    "D12" = 27.243215089760614
    "E12" = 21.74758906882867

    "B13" = 6.4456210257433533
    "C13" = 6.6863112203408237
    "D13" = 6.1533091811755591
    "E13" = 7.3800098108694145

    "B15" = 1.0644562102574333
    "C15" = 1.0668631122034082
    "D15" = 1.0615330918117374
    "E15" = 1.0738000981086941

    "B16" = 10.376554596361856
    "C16" = 15.538486585096859
    "D16" = 12.582062933265011
    "E16" = 7.1419521408909485
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
